# MP120_Transform.xlsx update
# - Drop the now-unused Sheet2 / Sheet3 tabs (workbook now holds just Sheet1)
# - Refresh the sensorCG sample values on Sheet1 (B2:D4) to the latest run
# - Re-select the full data range A1:D4 on Sheet1

$wb = $excel.ActiveWorkbook

# Remove Sheet2 and Sheet3, keeping only Sheet1
[void]$wb.Worksheets("Sheet2").Delete()
[void]$wb.Worksheets("Sheet3").Delete()

$ws = $wb.Worksheets("Sheet1")

# Updated numeric values for B2:D4
$ws.Range("B2").Value = 0.9780985493512234
$ws.Range("C2").Value = 0.18502749614949032
$ws.Range("D2").Value = -0.09533128251357295

$ws.Range("B3").Value = 0.1898478116266182
$ws.Range("C3").Value = -0.6052921505099105
$ws.Range("D3").Value = 0.7730324837622751

$ws.Range("B4").Value = 0.08532898790924687
$ws.Range("C4").Value = -0.7742003863340153
$ws.Range("D4").Value = -0.627162439582159

# Re-select the data range A1:D4 so that's the active selection on the sheet
[void]$ws.Range("A1:D4").Select()
